$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1421.6774
$ws.Range("I28").Value = 1589.7407
$ws.Range("J28").Value = 287.25
$ws.Range("K28").Value = 1589.7407
$ws.Range("L28").Value = 287.25
$ws.Range("M28").Value = -1104.7407
$ws.Range("N28").Value = -1257.25
$ws.Range("H33").Value = 250.76666
$ws.Range("J33").Value = 487.4
$ws.Range("L33").Value = 487.4
$ws.Range("N33").Value = -945.4
$ws.Range("H74").Value = 3314.7144
$ws.Range("I74").Value = 2950.75
$ws.Range("J74").Value = 3800
$ws.Range("K74").Value = 2950.75
$ws.Range("L74").Value = 3800
$ws.Range("M74").Value = -2014.75
$ws.Range("N74").Value = -5672
$ws.Range("H77").Value = 3314.7144
$ws.Range("I77").Value = 2950.75
$ws.Range("J77").Value = 3800
$ws.Range("K77").Value = 14753.75
$ws.Range("L77").Value = 19000
$ws.Range("M77").Value = -10073.75
$ws.Range("H96").Value = 1622.75
$ws.Range("I96").Value = 2134.125
$ws.Range("K96").Value = 6402.375
$ws.Range("M96").Value = -5029.375
$ws.Range("H101").Value = 933.3333
$ws.Range("I101").Value = 900
$ws.Range("K101").Value = 2700
$ws.Range("M101").Value = -1078
$ws.Range("H107").Value = 3066.0908
$ws.Range("I107").Value = 3296.1538
$ws.Range("J107").Value = 2733.7778
$ws.Range("K107").Value = 3296.1538
$ws.Range("L107").Value = 2733.7778
$ws.Range("M107").Value = -1376.1538
$ws.Range("N107").Value = -6573.7778
$ws.Range("H116").Value = 3836
$ws.Range("I116").Value = 2483
$ws.Range("K116").Value = 2483
$ws.Range("M116").Value = 959
$ws.Range("H121").Value = 1415
$ws.Range("J121").Value = 1393.75
$ws.Range("L121").Value = 4181.25
$ws.Range("N121").Value = -7675.25
$ws.Range("H132").Value = 6065907.5
$ws.Range("I132").Value = 8134737.5
$ws.Range("J132").Value = 7191.857
$ws.Range("K132").Value = 24404212.5
$ws.Range("L132").Value = 21575.571
$ws.Range("M132").Value = -24401682.5
$ws.Range("N132").Value = -26635.571
$ws.Range("H137").Value = 1446.025
$ws.Range("J137").Value = 1794.762
$ws.Range("L137").Value = 5384.286
$ws.Range("N137").Value = -10484.286
$ws.Range("H138").Value = 1637.96
$ws.Range("I138").Value = 1085.9166
$ws.Range("J138").Value = 1812.2894
$ws.Range("K138").Value = 3257.7498
$ws.Range("L138").Value = 5436.8682
$ws.Range("M138").Value = 1882.2502
$ws.Range("N138").Value = -15716.8682

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 656.38464
$ws.Range("I97").Value = 484.81818
$ws.Range("K97").Value = 484.81818
$ws.Range("M97").Value = 11.18182000000002
$ws.Range("H122").Value = 1999.4166
$ws.Range("I122").Value = 1397.9
$ws.Range("K122").Value = 4193.700000000001
$ws.Range("M122").Value = -1743.700000000001
$ws.Range("H132").Value = 2690.3635
$ws.Range("I132").Value = 2364.875
$ws.Range("J132").Value = 3558.3333
$ws.Range("K132").Value = 7094.625
$ws.Range("L132").Value = 10674.9999
$ws.Range("M132").Value = -4564.625
$ws.Range("N132").Value = -15734.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1152.9412
$ws.Range("I107").Value = 871.75
$ws.Range("K107").Value = 871.75
$ws.Range("M107").Value = 1048.25
$ws.Range("H134").Value = 4976.1787
$ws.Range("I134").Value = 935.0476
$ws.Range("J134").Value = 17099.572
$ws.Range("K134").Value = 2805.1428
$ws.Range("L134").Value = 51298.716
$ws.Range("M134").Value = -270.1428000000001
$ws.Range("N134").Value = -56368.716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1417.1212
$ws.Range("I31").Value = 1391.8667
$ws.Range("J31").Value = 1438.1666
$ws.Range("K31").Value = 1391.8667
$ws.Range("L31").Value = 1438.1666
$ws.Range("M31").Value = -1096.8667
$ws.Range("N31").Value = -2028.1666
$ws.Range("H32").Value = 3675
$ws.Range("I32").Value = 1566.6666
$ws.Range("K32").Value = 1566.6666
$ws.Range("M32").Value = -1250.6666
$ws.Range("H34").Value = 1417.1212
$ws.Range("I34").Value = 1391.8667
$ws.Range("J34").Value = 1438.1666
$ws.Range("K34").Value = 1391.8667
$ws.Range("L34").Value = 1438.1666
$ws.Range("M34").Value = -1189.8667
$ws.Range("N34").Value = -1842.1666
$ws.Range("H35").Value = 525
$ws.Range("I35").Value = 525
$ws.Range("K35").Value = 525
$ws.Range("M35").Value = -231
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H94").Value = 1123.1666
$ws.Range("I94").Value = 680
$ws.Range("K94").Value = 680
$ws.Range("M94").Value = -229
$ws.Range("H132").Value = 2569.1052
$ws.Range("J132").Value = 2832
$ws.Range("L132").Value = 8496
$ws.Range("N132").Value = -13556
$ws.Range("H134").Value = 18520340
$ws.Range("I134").Value = 1956.0476
$ws.Range("K134").Value = 5868.142800000001
$ws.Range("M134").Value = -3333.142800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1049.5883
$ws.Range("J122").Value = 1196.9231
$ws.Range("L122").Value = 10772.3079
$ws.Range("N122").Value = -15672.3079
$ws.Range("H131").Value = 22226168
$ws.Range("I131").Value = 125000344
$ws.Range("J131").Value = 4725.054
$ws.Range("K131").Value = 375001032
$ws.Range("L131").Value = 14175.162
$ws.Range("M131").Value = -374995992
$ws.Range("N131").Value = -24255.162
$ws.Range("H132").Value = 767.875
$ws.Range("I132").Value = 767.875
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6910.875
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4380.875
$ws.Range("N132").ClearContents()
$ws.Range("H139").Value = 1732.8
$ws.Range("I139").Value = 1811.2727
$ws.Range("J139").Value = 1600
$ws.Range("K139").Value = 5433.8181
$ws.Range("L139").Value = 4800
$ws.Range("M139").Value = -293.8181000000004
$ws.Range("N139").Value = -15080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4525.1113
$ws.Range("I80").Value = 3850
$ws.Range("J80").Value = 5065.2
$ws.Range("K80").Value = 3850
$ws.Range("L80").Value = 5065.2
$ws.Range("M80").Value = -2852
$ws.Range("N80").Value = -7061.2
$ws.Range("H83").Value = 4525.1113
$ws.Range("I83").Value = 3850
$ws.Range("J83").Value = 5065.2
$ws.Range("K83").Value = 19250
$ws.Range("L83").Value = 25326
$ws.Range("M83").Value = -14258
$ws.Range("N83").Value = -35310
$ws.Range("H122").Value = 1377.7838
$ws.Range("I122").Value = 1302.1428
$ws.Range("J122").Value = 1613.1111
$ws.Range("K122").Value = 3906.4284
$ws.Range("L122").Value = 4839.3333
$ws.Range("M122").Value = -1456.4284
$ws.Range("N122").Value = -9739.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 314.1111
$ws.Range("I55").Value = 288.2
$ws.Range("J55").Value = 346.5
$ws.Range("K55").Value = 288.2
$ws.Range("L55").Value = 346.5
$ws.Range("M55").Value = -115.2
$ws.Range("N55").Value = -692.5
$ws.Range("H82").Value = 2537.5
$ws.Range("I82").Value = 2500
$ws.Range("J82").Value = 2800
$ws.Range("K82").Value = 2500
$ws.Range("L82").Value = 2800
$ws.Range("M82").Value = -2139
$ws.Range("N82").Value = -3522
$ws.Range("H85").Value = 2537.5
$ws.Range("I85").Value = 2500
$ws.Range("J85").Value = 2800
$ws.Range("K85").Value = 2500
$ws.Range("L85").Value = 2800
$ws.Range("M85").Value = -1252
$ws.Range("N85").Value = -5296
$ws.Range("H93").Value = 1050
$ws.Range("I93").Value = 1000
$ws.Range("J93").Value = 1100
$ws.Range("K93").Value = 1000
$ws.Range("L93").Value = 1100
$ws.Range("M93").Value = 248
$ws.Range("N93").Value = -3596
$ws.Range("H138").Value = 34000
$ws.Range("J138").Value = 34000
$ws.Range("L138").Value = 34000
$ws.Range("N138").Value = -44280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 333335260
$ws.Range("I26").Value = 1000000000
$ws.Range("J26").Value = 2900
$ws.Range("K26").Value = 1000000000
$ws.Range("L26").Value = 2900
$ws.Range("M26").Value = -999999707
$ws.Range("N26").Value = -3486
$ws.Range("H28").Value = 41012.668
$ws.Range("J28").Value = 41012.668
$ws.Range("L28").Value = 41012.668
$ws.Range("N28").Value = -41708.668
$ws.Range("H29").Value = 800
$ws.Range("I29").Value = 800
$ws.Range("K29").Value = 800
$ws.Range("M29").Value = -510
$ws.Range("H31").Value = 3000
$ws.Range("J31").Value = 3000
$ws.Range("L31").Value = 3000
$ws.Range("N31").Value = -3696
$ws.Range("H39").Value = 1750
$ws.Range("I39").Value = 500
$ws.Range("J39").Value = 3000
$ws.Range("K39").Value = 500
$ws.Range("L39").Value = 3000
$ws.Range("M39").Value = -87
$ws.Range("N39").Value = -3826
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H108").Value = 27999.5
$ws.Range("J108").Value = 27999.5
$ws.Range("L108").Value = 27999.5
$ws.Range("N108").Value = -35679.5
$ws.Range("H122").Value = 12501002
$ws.Range("I122").Value = 14707011
$ws.Range("J122").Value = 286.66666
$ws.Range("K122").Value = 44121033
$ws.Range("L122").Value = 859.9999799999999
$ws.Range("M122").Value = -44118583
$ws.Range("N122").Value = -5759.99998
$ws.Range("H132").Value = 7102.4287
$ws.Range("I132").Value = 9954.5
$ws.Range("K132").Value = 29863.5
$ws.Range("M132").Value = -27333.5
$ws.Range("H136").Value = 1127.4138
$ws.Range("I136").Value = 1110.1
$ws.Range("K136").Value = 3330.3
$ws.Range("M136").Value = -780.2999999999997
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
